$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("J12")
Write-Host "J12 Interior.Color:" $c.Interior.Color
Write-Host "J12 Interior.ColorIndex:" $c.Interior.ColorIndex
Write-Host "J12 Interior.Pattern:" $c.Interior.Pattern

$c2 = $ws.Range("B12")
Write-Host "B12 Interior.Color:" $c2.Interior.Color
Write-Host "B12 Font.Color:" $c2.Font.Color

$c3 = $ws.Range("C4")
Write-Host "C4(blue) Interior.Color:" $c3.Interior.Color
